$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 6-9 (ECs/FAPs x MuSCs and Resolving-Mac x ECs/MuSCs pairs)
$ws.Range("A6:T9").EntireRow.Delete() | Out-Null

# Row 2: ECs -> Rtn4/Tnfrsf19 -> MuSCs (new TPM values)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Rtn4"
$ws.Range("C2").Value = "Tnfrsf19"
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 35.71508866666667
$ws.Range("H2").Value = 107.145266
$ws.Range("I2").Value = 0.1390302752364672
$ws.Range("J2").Value = 0.1390302752364672
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.769432666666667
$ws.Range("N2").Value = 5.308298
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 63.19544457969646
$ws.Range("R2").Value = 568.7590012172681
$ws.Range("S2").Value = 0.1390302752364672
$ws.Range("T2").Value = 0.1390302752364672

# Row 3: FAPs -> Rtn4/Tnfrsf19 -> MuSCs (new TPM values)
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Rtn4"
$ws.Range("C3").Value = "Tnfrsf19"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 54.09018966666667
$ws.Range("H3").Value = 162.270569
$ws.Range("I3").Value = 0.2105601368412127
$ws.Range("J3").Value = 0.2105601368412127
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.769432666666667
$ws.Range("N3").Value = 5.308298
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 95.70894854239577
$ws.Range("R3").Value = 861.380536881562
$ws.Range("S3").Value = 0.2105601368412127
$ws.Range("T3").Value = 0.2105601368412127

# Row 4: MuSCs -> Rtn4/Tnfrsf19 -> MuSCs (new TPM values)
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Rtn4"
$ws.Range("C4").Value = "Tnfrsf19"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 101.4529346666666
$ws.Range("H4").Value = 304.358804
$ws.Range("I4").Value = 0.3949319449238378
$ws.Range("J4").Value = 0.3949319449238378
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.769432666666667
$ws.Range("N4").Value = 5.308298
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 179.5141367283991
$ws.Range("R4").Value = 1615.627230555592
$ws.Range("S4").Value = 0.3949319449238378
$ws.Range("T4").Value = 0.3949319449238378

# Row 5: Resolving-Mac -> Rtn4/Tnfrsf19 -> MuSCs (new TPM values)
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("B5").Value = "Rtn4"
$ws.Range("C5").Value = "Tnfrsf19"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 65.628919
$ws.Range("H5").Value = 196.886757
$ws.Range("I5").Value = 0.2554776429984823
$ws.Range("J5").Value = 0.2554776429984823
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.769432666666667
$ws.Range("N5").Value = 5.308298
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 116.1259531566207
$ws.Range("R5").Value = 1045.133578409586
$ws.Range("S5").Value = 0.2554776429984823
$ws.Range("T5").Value = 0.2554776429984823
